$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

$ws.Range("M8").Value = "nan"
$ws.Range("O8").Value = "م.محمد عبدالله ،تيم الكرد"
